# Generate Report for Handoff
# Replaces the two "in-flight" source file rows (a .png and a .md) on each
# sheet with two new markdown files that have completed handoff, and removes
# the row that is no longer part of the report (row count goes 5 -> 4 on every
# sheet). Also refreshes the handoff timestamps / xlf target names to the new
# files' generated names.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 4 (bc2883ff...png) is dropped entirely; row 5 (.localization-config)
# shifts up into row 4.
$ov.Rows.Item(4).Delete()

# Rename the two still-present source files.
$ov.Range("A2").Value = "42970a04-655a-415e-a3ca-0937883cecb1.md"
$ov.Range("A3").Value = "5ca5d609-f787-4505-ad48-d3568810b619.md"

# Rebuild hyperlinks to match the new file names / row positions.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/e2e/42970a04-655a-415e-a3ca-0937883cecb1.md", $missing, $missing, "42970a04-655a-415e-a3ca-0937883cecb1.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/e2e/5ca5d609-f787-4505-ad48-d3568810b619.md", $missing, $missing, "5ca5d609-f787-4505-ad48-d3568810b619.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | Status | Latest Handoff File | ...
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 4 (bc2883ff...png, an IsDependency row) is dropped; row 5
# (.localization-config) shifts up into row 4.
$zh.Rows.Item(4).Delete()

$zh.Range("A2").Value = "42970a04-655a-415e-a3ca-0937883cecb1.md"
$zh.Range("C2").Value = "42970a04-655a-415e-a3ca-0937883cecb1.b6af4eb6627a9b2f27e3052a1bc102c1e22f5961.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-10 03:43:32"
$zh.Range("H2").Value = "Include"
$zh.Range("I2").ClearContents()

$zh.Range("A3").Value = "5ca5d609-f787-4505-ad48-d3568810b619.md"
$zh.Range("C3").Value = "5ca5d609-f787-4505-ad48-d3568810b619.85bd0e0ea2bba4d00b119158a99ef4546d1e7fe9.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-10 03:43:32"

# Row 4 is now the former ".localization-config" row; it no longer carries a
# "Dependency From" row (C4 cleared), and picks up the new row4 dates.
$zh.Range("B4").Value = "Not to be localized"
$zh.Range("C4").ClearContents()
$zh.Range("D4").Value = "0001-01-01 00:00:00"
$zh.Range("G4").Value = "0001-01-01 00:00:00"
$zh.Range("H4").Value = "Ignored"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/e2e/42970a04-655a-415e-a3ca-0937883cecb1.md", $missing, $missing, "42970a04-655a-415e-a3ca-0937883cecb1.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed2701e9eade08e514911b022d135d0698dfb908/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/42970a04-655a-415e-a3ca-0937883cecb1.b6af4eb6627a9b2f27e3052a1bc102c1e22f5961.zh-cn.xlf", $missing, $missing, "42970a04-655a-415e-a3ca-0937883cecb1.b6af4eb6627a9b2f27e3052a1bc102c1e22f5961.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/e2e/5ca5d609-f787-4505-ad48-d3568810b619.md", $missing, $missing, "5ca5d609-f787-4505-ad48-d3568810b619.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed2701e9eade08e514911b022d135d0698dfb908/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5ca5d609-f787-4505-ad48-d3568810b619.85bd0e0ea2bba4d00b119158a99ef4546d1e7fe9.zh-cn.xlf", $missing, $missing, "5ca5d609-f787-4505-ad48-d3568810b619.85bd0e0ea2bba4d00b119158a99ef4546d1e7fe9.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/.localization-config", $missing, $missing, ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": same shape as zh-cn, different target files/dates.
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Rows.Item(4).Delete()

$de.Range("A2").Value = "42970a04-655a-415e-a3ca-0937883cecb1.md"
$de.Range("C2").Value = "42970a04-655a-415e-a3ca-0937883cecb1.b6af4eb6627a9b2f27e3052a1bc102c1e22f5961.de-de.xlf"
$de.Range("D2").Value = "2016-03-10 03:43:35"
$de.Range("H2").Value = "Include"
$de.Range("I2").ClearContents()

$de.Range("A3").Value = "5ca5d609-f787-4505-ad48-d3568810b619.md"
$de.Range("C3").Value = "5ca5d609-f787-4505-ad48-d3568810b619.85bd0e0ea2bba4d00b119158a99ef4546d1e7fe9.de-de.xlf"
$de.Range("D3").Value = "2016-03-10 03:43:35"

$de.Range("B4").Value = "Not to be localized"
$de.Range("C4").ClearContents()
$de.Range("D4").Value = "0001-01-01 00:00:00"
$de.Range("G4").Value = "0001-01-01 00:00:00"
$de.Range("H4").Value = "Ignored"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/e2e/42970a04-655a-415e-a3ca-0937883cecb1.md", $missing, $missing, "42970a04-655a-415e-a3ca-0937883cecb1.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efd881b91a82455a9a8a3e24e46d07bebee0e6d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/42970a04-655a-415e-a3ca-0937883cecb1.b6af4eb6627a9b2f27e3052a1bc102c1e22f5961.de-de.xlf", $missing, $missing, "42970a04-655a-415e-a3ca-0937883cecb1.b6af4eb6627a9b2f27e3052a1bc102c1e22f5961.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/e2e/5ca5d609-f787-4505-ad48-d3568810b619.md", $missing, $missing, "5ca5d609-f787-4505-ad48-d3568810b619.md") | Out-Null
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efd881b91a82455a9a8a3e24e46d07bebee0e6d2/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5ca5d609-f787-4505-ad48-d3568810b619.85bd0e0ea2bba4d00b119158a99ef4546d1e7fe9.de-de.xlf", $missing, $missing, "5ca5d609-f787-4505-ad48-d3568810b619.85bd0e0ea2bba4d00b119158a99ef4546d1e7fe9.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/62a121aac7973910b8abc5f253d183002d2757a1/.localization-config", $missing, $missing, ".localization-config") | Out-Null

Write-Output "Report regenerated for handoff."
